$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (rows 2-6) ---
$ws.Range("A2").Value = "AD.SEC.001.FON.02"
$ws.Range("A3").Value = "AD.SEC.001.FON.01"
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A4").Value = "AD.SEC.001.FON.03"
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A5").Value = "RO.ACT"
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A6").Value = "AD.SEC.014.FON.01"

# --- Column D (rows 2-5) ---
$ws.Range("D2").Value = "AD.SEC.001.FON.02"
$ws.Range("D3").Value = "AD.SEC.001.FON.01"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").Value = "AD.SEC.001.FON.03"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").Value = "AD.SEC.014.FON.01"

# --- Old lookup block (rows 7-14, columns D:H) gets cleared ---
$ws.Range("D7:D8").ClearContents()
$ws.Range("E8:H14").ClearContents()

# --- New lookup block (rows 9-13, column D) ---
$ws.Range("D9").Value = "AD.SEC.001.FON.02"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").Value = "AD.SEC.001.FON.01"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").Value = "AD.SEC.001.FON.03"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").Value = "RO.ACT"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").Value = "AD.SEC.014.FON.01"

# --- Residual "Normal 18" cell style (unused, carried over from a copy/paste) ---
$st = $wb.Styles.Add("Normal 18")
$st.Font.Color = 0
$st.Font.Size = 10

# --- Selection ---
$null = $ws.Range("B12").Select()
